$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (time-to-something) and C (Energy) for rows 2-50
$data = @(
    @(2, 3.869511384177685, 1.707596222450402),
    @(3, 4.223292528482742, 3.117793797679763),
    @(4, 6.496134057827351, 5.218938197714804),
    @(5, 11.3329550434594, 7.032728163196263),
    @(6, 11.89843846762003, 8.486844331433653),
    @(7, 12.8945500185167, 10.20019337189851),
    @(8, 15.19166175243297, 11.94443904954075),
    @(9, 16.22120194155941, 13.89398399358445),
    @(10, 16.87591993044459, 16.28353721578357),
    @(11, 19.42619500186913, 17.99682078664295),
    @(12, 19.68594413535853, 19.80095201177858),
    @(13, 21.62985158346695, 21.49060152087361),
    @(14, 23.12136261958922, 23.3953018171421),
    @(15, 25.72256442315136, 24.97892517350413),
    @(16, 31.01143679269622, 26.47845004308031),
    @(17, 35.57861987815435, 28.34693683098946),
    @(18, 35.92067933791748, 29.981913006492),
    @(19, 36.19852036439411, 31.52742380049328),
    @(20, 36.38379672189063, 33.20120034966937),
    @(21, 40.14848560852405, 34.8496464834083),
    @(22, 42.08506123289358, 36.45304452571201),
    @(23, 43.77644943764468, 38.02218668434978),
    @(24, 48.75728942615456, 39.85955706143929),
    @(25, 50.53469804415509, 41.4990801615341),
    @(26, 53.59037083383997, 43.62629840698037),
    @(27, 54.5304760030787, 45.41338134816947),
    @(28, 57.26900158737699, 47.71939195090843),
    @(29, 60.658345521985, 49.05176400210095),
    @(30, 61.2488003530978, 50.97934376722848),
    @(31, 62.36075767500014, 52.38356825665204),
    @(32, 62.72145846868706, 54.02176661508813),
    @(33, 68.67692755083833, 55.69332249942369),
    @(34, 68.98567947486916, 57.31135354242321),
    @(35, 71.08342223268404, 59.44436527019901),
    @(36, 71.48431759889107, 61.00311326219551),
    @(37, 73.99614573686954, 63.19461093187064),
    @(38, 74.32281553909783, 65.11716001809319),
    @(39, 75.23699641228237, 67.42830037898507),
    @(40, 76.60464258122951, 68.66591502085771),
    @(41, 77.91010054956246, 71.47232034258454),
    @(42, 80.06493197273107, 73.180698696758),
    @(43, 82.05554453036711, 74.93723382676154),
    @(44, 82.56367298388226, 76.86729913136263),
    @(45, 85.98439130149596, 78.4603629533339),
    @(46, 87.66164533876054, 80.59731627271253),
    @(47, 89.06810087745794, 82.10491810817611),
    @(48, 91.80822857006503, 83.85674076038288),
    @(49, 93.17572466200814, 85.44711790686908),
    @(50, 97.40526171535366, 86.92520729805804)
)

foreach ($row in $data) {
    $r = $row[0]
    $bVal = $row[1]
    $cVal = $row[2]
    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = $cVal
}
